$d = $word.ActiveDocument

# 1. Update the DATE line: 2024-04-19 -> 2024-04-24
$d.Content.Find.Execute("2024-04-19", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-04-24", 2) | Out-Null

# 2. "Six key demographic components" -> "Five key demographic components"
$d.Content.Find.Execute("Six key demographic components", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Five key demographic components", 2) | Out-Null

# 3. Remove the duplicate "Indigenous Identity" bullet group (heading + its two
#    explanatory bullets that follow "(Distinctions Based) Indigenous Identity"),
#    identified by its distinctive first bullet text, BEFORE renaming the earlier
#    heading so the search text stays unambiguous.
$startPara = $null
$endPara = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text
    if ($text -eq "Because current DIP datasets are limited in their capability to provide Indigenous demographics, we also rolled up the Distinctions based demography to a singular Indigenous/Non Indigenous option.`r") {
        $startPara = $d.Paragraphs.Item($i - 1)
        $endPara = $d.Paragraphs.Item($i + 1)
        break
    }
}
if ($startPara -ne $null) {
    $deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $deleteRange.Delete()
}

# 4. Rename the "(Distinctions Based) Indigenous Identity" heading to "Indigenous Identity"
$d.Content.Find.Execute("(Distinctions Based) Indigenous Identity", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Indigenous Identity", 2) | Out-Null

Write-Output "done"
